# Fix aspect ratio issues: update Res_Bounding_Box (A) and SK_Bounding_Box (B)
# coordinate values on the "Mapping" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -87.9168
$ws.Range("B2").Value = -87.8859

$ws.Range("A3").Value = 30.3022
$ws.Range("B3").Value = 30.329

$ws.Range("A4").Value = -87.7052
$ws.Range("B4").Value = -87.7362

$ws.Range("A5").Value = 30.4849
$ws.Range("B5").Value = 30.4582
